# Panel (d) box correction on slide 1 ("TextBox 61", shape id 62):
#  1. Fix a typo in its title ("ECFDs" -> "ECDFs").
#  2. Append a blank line plus an explanatory note about the grey
#     reference CDF.
#  3. Resize/reposition the legend box (done last: the box has
#     <a:spAutoFit/>, so resizing before the text edits would just get
#     overridden by the auto-fit recalculation triggered by the new
#     paragraphs).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(13)

$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- 1) Fix "ECFDs" -> "ECDFs" typo in the title run ---------------------
$oldTitle = "(SUBJECTIVE) CATEGORIES OF SIMILARITY DEGREES BETWEEN ECFDs"
$newTitle = "(SUBJECTIVE) CATEGORIES OF SIMILARITY DEGREES BETWEEN ECDFs"
$fullText = $tr.Text
$titleIdx = $fullText.IndexOf($oldTitle)
$titleRange = $tr.Characters($titleIdx + 1, $oldTitle.Length)
$titleRange.Text = $newTitle

# --- 2) Append a blank paragraph + explanatory note -----------------------
$beforeLen = $tr.Length
$noteText = "CDF in grey -> reference distribution computed from rain gauge observations."
$tr.InsertAfter("`r`r" + $noteText) | Out-Null

# Blank paragraph: same pink used by the "Very poor" entry above it.
$blankRange = $tr.Characters($beforeLen + 1, 1)
$blankRange.Font.Color.RGB = 0x6600FF   # BGR encoding of RGB FF0066

# Note paragraph: grey text (bg1 at 50% luminance -> mid-grey 808080).
$noteRange = $tr.Characters($beforeLen + 3, $noteText.Length)
$noteRange.Font.Color.RGB = 0x808080

# --- 3) Reposition / resize the box ---------------------------------------
# Target EMU values: off x=1950157 y=4434622 ; ext cx=3621040 cy=1261884
# Shape.Left/Top/Width/Height round-trip through single-precision points
# in this object model, so the literals below are the doubles whose
# float32 rounding reproduces those exact EMU values on save.
$sh.Left   = 153.55567179133857
$sh.Top    = 349.1828461456693
$sh.Width  = 285.1212616425197
$sh.Height = 99.36094288188977
